$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark after "Allow admin to create new batches"
# ---------------------------------------------------------------------------
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
}

# ---------------------------------------------------------------------------
# 2) Merge the three runs "Varchar(" + "100" + ")" into a single run "Varchar(100)"
#    (the Email column's data type cell in the 4th table)
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(4)
$cell = $t.Cell(2, 3)
$p = $cell.Range.Paragraphs.Item(1)
$pStart = $p.Range.Start
# run1 = "Varchar(" (8 chars), run2 = "100" (3 chars), run3 = ")" (1 char)
$run3 = $d.Range($pStart + 11, $pStart + 12)
$run3.Delete()
$run2 = $d.Range($pStart + 8, $pStart + 11)
$run2.Delete()
$run1 = $d.Range($pStart, $pStart + 8)
$run1.Text = "Varchar(100)"

# ---------------------------------------------------------------------------
# 3) Append the new "Server Pages" documentation section at the end of the
#    document (after the last table, before the two trailing empty paragraphs)
# ---------------------------------------------------------------------------
$secondEmpty = $d.Paragraphs.Last
$firstEmpty = $secondEmpty.Previous()

$xmlBlockA = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>SeverPages</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Server pages for the site serve one of </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">three </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">main purposes: displaying a UI, </w:t>
            </w:r>
            <w:r>
              <w:t>providing data in JSON format for client side processing</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> and processing submissions from </w:t>
            </w:r>
            <w:r>
              <w:t>users</w:t>
            </w:r>
            <w:r>
              <w:t>. Pages and their roles are listed below.</w:t>
            </w:r>
          </w:p>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insAtA = $d.Range($firstEmpty.Range.Start, $firstEmpty.Range.Start)
$insAtA.InsertXML($xmlBlockA)

# InsertXML always leaves one extra blank paragraph behind (to avoid merging
# its last paragraph into $firstEmpty); remove that spare blank paragraph so
# $firstEmpty remains the single blank separator paragraph.
$secondEmptyA = $d.Paragraphs.Last
$firstEmptyA = $secondEmptyA.Previous()
$spareBlankA = $firstEmptyA.Previous()
$delRangeA = $d.Range($spareBlankA.Range.Start, $firstEmptyA.Range.Start)
$delRangeA.Delete()

# ---------------------------------------------------------------------------
# 4) Insert the page-by-page descriptions (Index / Batches / sample JSON /
#    Register) between the (preserved) blank paragraph and the final blank
#    paragraph, plus the new _GoBack bookmark after the Register paragraph.
# ---------------------------------------------------------------------------
$secondEmptyB0 = $d.Paragraphs.Last
$firstEmptyB0 = $secondEmptyB0.Previous()

$xmlBlockB = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>Index: entry page to the site. Displays the main public UI which is the batch grid.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Batches: retrieves </w:t>
            </w:r>
            <w:r>
              <w:t>batch information for pending and ongoing batches, as well as batch-timeslot pairings. Data will reflect below.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:i/>
                <w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/>
              </w:rPr>
              <w:t>{id:</w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/>
              </w:rPr>
              <w:t>1,status</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/>
              </w:rPr>
              <w:t>:1,summary:"short description",courseid:1,coursename:"course 1",availseats:10,maxseats:20,startdate:'2017-06-20',enddate:'2017-08-20',teacher:'Subra'}</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/>
              </w:rPr>
              <w:t>{timeslot:</w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/>
              </w:rPr>
              <w:t>0,batchId</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/>
              </w:rPr>
              <w:t>:1}</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>Register: receives and processes student registration.</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insAtB = $d.Range($firstEmptyB0.Range.Start, $firstEmptyB0.Range.Start)
$insAtB.InsertXML($xmlBlockB)

# Remove the spare blank paragraph InsertXML left behind so the original
# trailing blank paragraph remains the sole separator before the sectPr.
$secondEmptyB = $d.Paragraphs.Last
$firstEmptyB = $secondEmptyB.Previous()
$spareBlankB = $firstEmptyB.Previous()
$delRangeB = $d.Range($spareBlankB.Range.Start, $firstEmptyB.Range.Start)
$delRangeB.Delete()

Write-Output "Done"
